# Refresh the cryptos price/volume snapshot (GitHub Actions scrape update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($range, $value) {
    # Price column holds numeric-looking text (e.g. "1.003", "26.061.75");
    # force text format so Excel does not reinterpret it as a number/date,
    # then restore the cell to its original (default) style.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextCell $ws.Range("D2") "26.061.75"
$ws.Range("E2").Value = "  -0.18%  "
Set-TextCell $ws.Range("D3") "1.646.83"
$ws.Range("E3").Value = "  +0.03%  "
Set-TextCell $ws.Range("D5") "218.21"
$ws.Range("E5").Value = "  +0.31%  "
Set-TextCell $ws.Range("D6") "0.5194"
$ws.Range("E6").Value = "  -0.12%  "
Set-TextCell $ws.Range("D7") "1.003"
$ws.Range("E7").Value = "  -0.15%  "
Set-TextCell $ws.Range("D8") "0.2621"
$ws.Range("E8").Value = "  +0.17%  "
Set-TextCell $ws.Range("D9") "0.06297"
$ws.Range("E9").Value = "  +0.27%  "
Set-TextCell $ws.Range("D10") "20.27"
$ws.Range("E10").Value = "  -0.84%  "
Set-TextCell $ws.Range("D11") "0.07684"
$ws.Range("E11").Value = "  -0.89%  "
Set-TextCell $ws.Range("D12") "4.592"
$ws.Range("E12").Value = "  +2.68%  "
Set-TextCell $ws.Range("D13") "1.649.58"
$ws.Range("E13").Value = "  -1.06%  "
Set-TextCell $ws.Range("D14") "1.872.81"
$ws.Range("E14").Value = "  +0.09%  "
Set-TextCell $ws.Range("D15") "0.5567"
$ws.Range("E15").Value = "  -0.38%  "
Set-TextCell $ws.Range("D16") "0.0₅8092"
$ws.Range("E16").Value = "  +1.10%  "
Set-TextCell $ws.Range("D17") "65.02"
$ws.Range("E17").Value = "  +0.33%  "
Set-TextCell $ws.Range("D18") "26.037.37"
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("E19").Value = "  -0.19%  "
Set-TextCell $ws.Range("D20") "4.596"
$ws.Range("E20").Value = "  -1.13%  "
Set-TextCell $ws.Range("D21") "192.88"
$ws.Range("E21").Value = "  +0.30%  "
Set-TextCell $ws.Range("D22") "10.42"
$ws.Range("E22").Value = "  +3.04%  "
Set-TextCell $ws.Range("D23") "5.909"
$ws.Range("E23").Value = "  -0.69%  "
Set-TextCell $ws.Range("D25") "144.62"
$ws.Range("E25").Value = "  -1.07%  "
Set-TextCell $ws.Range("D26") "0.1178"
$ws.Range("E26").Value = "  -1.81%  "
Set-TextCell $ws.Range("D27") "7.181"
$ws.Range("E27").Value = "  +0.10%  "
Set-TextCell $ws.Range("D28") "15.82"
$ws.Range("E28").Value = "  -0.59%  "
$ws.Range("E29").Value = "  +2.29%  "
Set-TextCell $ws.Range("D30") "0.05347"
$ws.Range("E30").Value = "  -4.75%  "
Set-TextCell $ws.Range("D31") "1.266"
$ws.Range("E31").Value = "  +0.18%  "
Set-TextCell $ws.Range("D32") "3.445"
Set-TextCell $ws.Range("D33") "3.323"
$ws.Range("E33").Value = "  -0.63%  "
Set-TextCell $ws.Range("D34") "1.546"
$ws.Range("E34").Value = "  -3.01%  "
Set-TextCell $ws.Range("D35") "2.417"
$ws.Range("E35").Value = "  +0.27%  "
$ws.Range("E36").Value = "  -0.26%  "
Set-TextCell $ws.Range("D37") "0.9405"
$ws.Range("E37").Value = "  +0.44%  "
Set-TextCell $ws.Range("D38") "0.5574"
$ws.Range("E38").Value = "  -1.68%  "
Set-TextCell $ws.Range("D39") "0.01573"
$ws.Range("E39").Value = "  -0.62%  "
Set-TextCell $ws.Range("D40") "5.762"
$ws.Range("E40").Value = "  -3.51%  "
Set-TextCell $ws.Range("D42") "1.025.25"
$ws.Range("E42").Value = "  -2.37%  "
Set-TextCell $ws.Range("D43") "0.8233"
$ws.Range("E43").Value = "  -2.14%  "
Set-TextCell $ws.Range("D44") "100.71"
$ws.Range("E44").Value = "  -1.54%  "
Set-TextCell $ws.Range("D45") "1.782.40"
$ws.Range("E45").Value = "  +0.01%  "
$ws.Range("E46").Value = "  +2.38%  "
Set-TextCell $ws.Range("D47") "57.17"
$ws.Range("E47").Value = "  -0.15%  "
$ws.Range("E48").Value = "  -0.55%  "
Set-TextCell $ws.Range("D49") "0.4310"
$ws.Range("E49").Value = "  -0.56%  "
Set-TextCell $ws.Range("D50") "7.853"
$ws.Range("E50").Value = "  -0.61%  "
Set-TextCell $ws.Range("D51") "0.05100"
$ws.Range("E51").Value = "  -4.42%  "
